# Auto-generated script to update cryptos sheet values per commit diff
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($ws, $ref, $val) {
    $c = $ws.Range($ref)
    $c.NumberFormat = "@"
    $c.Value = $val
    $c.ClearFormats()
}

Set-TextCell $ws 'D2' '97.433.37'
Set-TextCell $ws 'E2' '  -1.11%  '
Set-TextCell $ws 'D3' '3.352.24'
Set-TextCell $ws 'E3' '  -0.28%  '
Set-TextCell $ws 'D4' '0.999'
Set-TextCell $ws 'E4' '  -0.07%  '
Set-TextCell $ws 'D5' '249.36'
Set-TextCell $ws 'E5' '  -3.00%  '
Set-TextCell $ws 'D6' '657.49'
Set-TextCell $ws 'E6' '  -1.05%  '
Set-TextCell $ws 'D7' '1.39'
Set-TextCell $ws 'E7' '  -9.79%  '
Set-TextCell $ws 'D8' '0.418'
Set-TextCell $ws 'E8' '  -11.64%  '
Set-TextCell $ws 'D9' '0.999'
Set-TextCell $ws 'E9' '  +0.01%  '
Set-TextCell $ws 'D10' '1.02'
Set-TextCell $ws 'E10' '  -4.90%  '
Set-TextCell $ws 'D11' '3.348.60'
Set-TextCell $ws 'E11' '  -0.31%  '
Set-TextCell $ws 'D12' '0.209'
Set-TextCell $ws 'E12' '  -3.06%  '
Set-TextCell $ws 'D13' '40.54'
Set-TextCell $ws 'E13' '  -4.33%  '
Set-TextCell $ws 'D14' '97.148.26'
Set-TextCell $ws 'E14' '  -1.28%  '
Set-TextCell $ws 'D15' '6.09'
Set-TextCell $ws 'E15' '  +6.82%  '
Set-TextCell $ws 'D16' '0.0000255'
Set-TextCell $ws 'E16' '  -7.15%  '
Set-TextCell $ws 'D17' '3.978.93'
Set-TextCell $ws 'E17' '  -0.32%  '
Set-TextCell $ws 'D18' '8.55'
Set-TextCell $ws 'E18' '  +11.02%  '
Set-TextCell $ws 'D19' '3.345.74'
Set-TextCell $ws 'E19' '  -0.31%  '
Set-TextCell $ws 'D20' '0.550'
Set-TextCell $ws 'E20' '  +26.18%  '
Set-TextCell $ws 'D21' '16.85'
Set-TextCell $ws 'E21' '  +0.44%  '
Set-TextCell $ws 'D22' '10.79'
Set-TextCell $ws 'E22' '  +2.09%  '
Set-TextCell $ws 'D23' '498.17'
Set-TextCell $ws 'E23' '  -6.36%  '
Set-TextCell $ws 'D24' '3.35'
Set-TextCell $ws 'E24' '  -6.75%  '
Set-TextCell $ws 'D25' '0.0000200'
Set-TextCell $ws 'E25' '  -8.30%  '
Set-TextCell $ws 'D26' '6.21'
Set-TextCell $ws 'E26' '  -0.25%  '
Set-TextCell $ws 'D27' '94.09'
Set-TextCell $ws 'E27' '  -8.40%  '
Set-TextCell $ws 'D28' '12.12'
Set-TextCell $ws 'E28' '  -3.68%  '
Set-TextCell $ws 'D29' '3.548.86'
Set-TextCell $ws 'E29' '  +0.31%  '
Set-TextCell $ws 'D30' '0.147'
Set-TextCell $ws 'E30' '  -1.19%  '
Set-TextCell $ws 'B31' 'Dai'
Set-TextCell $ws 'C31' 'https://coinranking.com/coin/MoTuySvg7+dai-dai'
Set-TextCell $ws 'D31' '0.996'
Set-TextCell $ws 'E31' '  -0.43%  '
Set-TextCell $ws 'B32' 'InternetComputer(DFINITY)'
Set-TextCell $ws 'C32' 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp'
Set-TextCell $ws 'D32' '11.04'
Set-TextCell $ws 'D33' '0.190'
Set-TextCell $ws 'E33' '  -0.21%  '
Set-TextCell $ws 'D34' '2.53'
Set-TextCell $ws 'E34' '  +19.52%  '
Set-TextCell $ws 'E35' '  -0.14%  '
Set-TextCell $ws 'D36' '0.551'
Set-TextCell $ws 'E36' '  +2.42%  '
Set-TextCell $ws 'D37' '28.49'
Set-TextCell $ws 'E37' '  -3.47%  '
Set-TextCell $ws 'D38' '7.71'
Set-TextCell $ws 'E38' '  -1.03%  '
Set-TextCell $ws 'E39' '  +10.26%  '
Set-TextCell $ws 'D40' '521.95'
Set-TextCell $ws 'E40' '  -0.85%  '
Set-TextCell $ws 'E41' '  -0.01%  '
Set-TextCell $ws 'D42' '0.149'
Set-TextCell $ws 'E42' '  -5.40%  '
Set-TextCell $ws 'D43' '24.66'
Set-TextCell $ws 'E43' '  -0.18%  '
Set-TextCell $ws 'D44' '8.82'
Set-TextCell $ws 'E44' '  +13.47%  '
Set-TextCell $ws 'D45' '0.842'
Set-TextCell $ws 'E45' '  +2.10%  '
Set-TextCell $ws 'D46' '0.0420'
Set-TextCell $ws 'E46' '  -3.37%  '
Set-TextCell $ws 'D47' '3.66'
Set-TextCell $ws 'E47' '  -6.22%  '
Set-TextCell $ws 'D48' '5.58'
Set-TextCell $ws 'E48' '  +8.97%  '
Set-TextCell $ws 'D49' '1.65'
Set-TextCell $ws 'E49' '  +7.24%  '
Set-TextCell $ws 'D50' '53.43'
Set-TextCell $ws 'E50' '  +4.62%  '
Set-TextCell $ws 'D51' '3.15'
Set-TextCell $ws 'E51' '  -8.21%  '
